$d = $word.ActiveDocument

# 1) "Averages:" -> "Weighted Averages:"
$d.Content.Find.Execute("Averages:", $true, $false, $false, $false, $false, $true, 1, $false, "Weighted Averages:", 2)

# 2) "24.41" -> "24.34"
$d.Content.Find.Execute("24.41", $true, $false, $false, $false, $false, $true, 1, $false, "24.34", 2)

# 3) "74.26" -> "" (clear the cell, but keep the run/text element)
$r1 = $d.Content
$r1.Find.Execute("74.26", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$r1.Text = ""

# 4) "1.33" -> "" (clear the cell, but keep the run/text element)
$r2 = $d.Content
$r2.Find.Execute("1.33", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$r2.Text = ""

# 5) "36.24" -> "34.98"
$d.Content.Find.Execute("36.24", $true, $false, $false, $false, $false, $true, 1, $false, "34.98", 2)
